# Scheduled market-data refresh: updates currentAveragePrice* / LevePrice* /
# LeveProfit* columns (H:N) across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and
# WVR Leve-profit tables with freshly pulled prices.
#
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 664.8570999999999  # H28 was 667.4761999999999
$ws.Cells.Item(28, 9).Value = 639.2778  # I28 was 665.41174
$ws.Cells.Item(28, 10).Value = 818.3333  # J28 was 676.25
$ws.Cells.Item(28, 11).Value = 639.2778  # K28 was 665.41174
$ws.Cells.Item(28, 12).Value = 818.3333  # L28 was 676.25
$ws.Cells.Item(28, 13).Value = -154.2778  # M28 was -180.41174
$ws.Cells.Item(28, 14).Value = -1788.3333  # N28 was -1646.25
$ws.Cells.Item(33, 8).Value = 246.38889  # H33 was 259.05884
$ws.Cells.Item(33, 9).Value = 162.86667  # I33 was 172.28572
$ws.Cells.Item(33, 11).Value = 162.86667  # K33 was 172.28572
$ws.Cells.Item(33, 13).Value = 66.13333  # M33 was 56.71428
$ws.Cells.Item(137, 8).Value = 3351.587  # H137 was 5560
$ws.Cells.Item(137, 9).Value = 2779.325  # I137 was 0
$ws.Cells.Item(137, 10).Value = 7166.6665  # J137 was 5560
$ws.Cells.Item(137, 11).Value = 8337.974999999999  # K137 was 0
$ws.Cells.Item(137, 12).Value = 21499.9995  # L137 was 16680
$ws.Cells.Item(137, 13).Value = -5787.974999999999  # M137 (new cell)
$ws.Cells.Item(137, 14).Value = -26599.9995  # N137 was -21780
$ws.Cells.Item(141, 8).Value = 2400.92  # H141 was 1812.6364
$ws.Cells.Item(141, 9).Value = 2060.1365  # I141 was 1722.4062
$ws.Cells.Item(141, 10).Value = 4900  # J141 was 4700
$ws.Cells.Item(141, 11).Value = 6180.4095  # K141 was 5167.2186
$ws.Cells.Item(141, 12).Value = 14700  # L141 was 14100
$ws.Cells.Item(141, 13).Value = -1000.4095  # M141 was 12.78139999999985
$ws.Cells.Item(141, 14).Value = -25060  # N141 was -24460
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 674.41174  # H2 was 513.7292
$ws.Cells.Item(2, 9).Value = 651.44446  # I2 was 487.64102
$ws.Cells.Item(2, 10).Value = 763  # J2 was 626.7778
$ws.Cells.Item(2, 11).Value = 651.44446  # K2 was 487.64102
$ws.Cells.Item(2, 12).Value = 763  # L2 was 626.7778
$ws.Cells.Item(2, 13).Value = -538.44446  # M2 was -374.64102
$ws.Cells.Item(2, 14).Value = -989  # N2 was -852.7778
$ws.Cells.Item(32, 8).Value = 6124.82  # H32 was 8441.111999999999
$ws.Cells.Item(32, 9).Value = 2990.8733  # I32 was 5410.379
$ws.Cells.Item(32, 10).Value = 13797.586  # J32 was 14692
$ws.Cells.Item(32, 11).Value = 2990.8733  # K32 was 5410.379
$ws.Cells.Item(32, 12).Value = 13797.586  # L32 was 14692
$ws.Cells.Item(32, 13).Value = -2703.8733  # M32 was -5123.379
$ws.Cells.Item(32, 14).Value = -14371.586  # N32 was -15266
$ws.Cells.Item(45, 8).Value = 2224.7273  # H45 was 2327.2
$ws.Cells.Item(45, 9).Value = 1170.6666  # I45 was 1156
$ws.Cells.Item(45, 11).Value = 1170.6666  # K45 was 1156
$ws.Cells.Item(45, 13).Value = -793.6666  # M45 was -779
$ws.Cells.Item(74, 8).Value = 1325.6136  # H74 was 1459.0264
$ws.Cells.Item(74, 9).Value = 964.64105  # I74 was 1052.6364
$ws.Cells.Item(74, 11).Value = 964.64105  # K74 was 1052.6364
$ws.Cells.Item(74, 13).Value = -90.64104999999995  # M74 was -178.6364000000001
$ws.Cells.Item(77, 8).Value = 1325.6136  # H77 was 1459.0264
$ws.Cells.Item(77, 9).Value = 964.64105  # I77 was 1052.6364
$ws.Cells.Item(77, 11).Value = 4823.20525  # K77 was 5263.182000000001
$ws.Cells.Item(77, 13).Value = -455.20525  # M77 was -895.1820000000007
$ws.Cells.Item(88, 8).Value = 22224688  # H88 was 11114844
$ws.Cells.Item(88, 10).Value = 0  # J88 was 5000
$ws.Cells.Item(88, 12).Value = 0  # L88 was 5000
$ws.Cells.Item(88, 14).ClearContents()  # N88 was -5812
$ws.Cells.Item(91, 8).Value = 22224688  # H91 was 11114844
$ws.Cells.Item(91, 10).Value = 0  # J91 was 5000
$ws.Cells.Item(91, 12).Value = 0  # L91 was 5000
$ws.Cells.Item(91, 14).ClearContents()  # N91 was -7808
$ws.Cells.Item(116, 8).Value = 674.41174  # H116 was 513.7292
$ws.Cells.Item(116, 9).Value = 651.44446  # I116 was 487.64102
$ws.Cells.Item(116, 10).Value = 763  # J116 was 626.7778
$ws.Cells.Item(116, 11).Value = 651.44446  # K116 was 487.64102
$ws.Cells.Item(116, 12).Value = 763  # L116 was 626.7778
$ws.Cells.Item(116, 13).Value = 1642.55554  # M116 was 1806.35898
$ws.Cells.Item(116, 14).Value = -5351  # N116 was -5214.7778
$ws.Cells.Item(132, 8).Value = 3283.4443  # H132 was 3864.9722
$ws.Cells.Item(132, 9).Value = 2619.027  # I132 was 3183.24
$ws.Cells.Item(132, 10).Value = 4729.5293  # J132 was 5414.364
$ws.Cells.Item(132, 11).Value = 7857.081  # K132 was 9549.719999999999
$ws.Cells.Item(132, 12).Value = 14188.5879  # L132 was 16243.092
$ws.Cells.Item(132, 13).Value = -5327.081  # M132 was -7019.719999999999
$ws.Cells.Item(132, 14).Value = -19248.5879  # N132 was -21303.092
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 674.41174  # H3 was 513.7292
$ws.Cells.Item(3, 9).Value = 651.44446  # I3 was 487.64102
$ws.Cells.Item(3, 10).Value = 763  # J3 was 626.7778
$ws.Cells.Item(3, 11).Value = 651.44446  # K3 was 487.64102
$ws.Cells.Item(3, 12).Value = 763  # L3 was 626.7778
$ws.Cells.Item(3, 13).Value = -537.44446  # M3 was -373.64102
$ws.Cells.Item(3, 14).Value = -991  # N3 was -854.7778
$ws.Cells.Item(22, 8).Value = 1028.5714  # H22 was 342
$ws.Cells.Item(22, 9).Value = 480  # I22 was 198.83333
$ws.Cells.Item(22, 10).Value = 2400  # J22 was 1201
$ws.Cells.Item(22, 11).Value = 480  # K22 was 198.83333
$ws.Cells.Item(22, 12).Value = 2400  # L22 was 1201
$ws.Cells.Item(22, 13).Value = -307  # M22 was -25.83332999999999
$ws.Cells.Item(22, 14).Value = -2746  # N22 was -1547
$ws.Cells.Item(105, 8).Value = 2939  # H105 was 3010
$ws.Cells.Item(105, 9).Value = 2939  # I105 was 3010
$ws.Cells.Item(105, 11).Value = 2939  # K105 was 3010
$ws.Cells.Item(105, 13).Value = -1192  # M105 was -1263
$ws.Cells.Item(134, 8).Value = 2295.9727  # H134 was 2407.836
$ws.Cells.Item(134, 9).Value = 1379  # I134 was 1456.4386
$ws.Cells.Item(134, 10).Value = 7464.364  # J134 was 7830.8
$ws.Cells.Item(134, 11).Value = 4137  # K134 was 4369.3158
$ws.Cells.Item(134, 12).Value = 22393.092  # L134 was 23492.4
$ws.Cells.Item(134, 13).Value = -1602  # M134 was -1834.3158
$ws.Cells.Item(134, 14).Value = -27463.092  # N134 was -28562.4
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2460.0208  # H31 was 2770.1191
$ws.Cells.Item(31, 9).Value = 900.67645  # I31 was 1031.6786
$ws.Cells.Item(31, 11).Value = 900.67645  # K31 was 1031.6786
$ws.Cells.Item(31, 13).Value = -605.67645  # M31 was -736.6786
$ws.Cells.Item(34, 8).Value = 2460.0208  # H34 was 2770.1191
$ws.Cells.Item(34, 9).Value = 900.67645  # I34 was 1031.6786
$ws.Cells.Item(34, 11).Value = 900.67645  # K34 was 1031.6786
$ws.Cells.Item(34, 13).Value = -698.67645  # M34 was -829.6786
$ws.Cells.Item(134, 8).Value = 4374.4863  # H134 was 4585.8857
$ws.Cells.Item(134, 9).Value = 4718.154  # I134 was 4890.88
$ws.Cells.Item(134, 10).Value = 3562.182  # J134 was 3823.4
$ws.Cells.Item(134, 11).Value = 14154.462  # K134 was 14672.64
$ws.Cells.Item(134, 12).Value = 10686.546  # L134 was 11470.2
$ws.Cells.Item(134, 13).Value = -11619.462  # M134 was -12137.64
$ws.Cells.Item(134, 14).Value = -15756.546  # N134 was -16540.2
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(42, 8).Value = 150000  # H42 was 0
$ws.Cells.Item(42, 9).Value = 150000  # I42 was 0
$ws.Cells.Item(42, 11).Value = 450000  # K42 was 0
$ws.Cells.Item(42, 13).Value = -449466  # M42 (new cell)
$ws.Cells.Item(113, 8).Value = 594.95776  # H113 was 634.88135
$ws.Cells.Item(113, 9).Value = 590.1731  # I113 was 638.85364
$ws.Cells.Item(113, 10).Value = 608.0526  # J113 was 625.8333
$ws.Cells.Item(113, 11).Value = 1770.5193  # K113 was 1916.56092
$ws.Cells.Item(113, 12).Value = 1824.1578  # L113 was 1877.4999
$ws.Cells.Item(113, 13).Value = 399.4807000000001  # M113 was 253.4390799999999
$ws.Cells.Item(113, 14).Value = -6164.1578  # N113 was -6217.4999
$ws.Cells.Item(130, 8).Value = 2592.7334  # H130 was 2842.4167
$ws.Cells.Item(130, 9).Value = 1000  # I130 was 1124
$ws.Cells.Item(130, 10).Value = 3986.375  # J130 was 4069.8572
$ws.Cells.Item(130, 11).Value = 3000  # K130 was 3372
$ws.Cells.Item(130, 12).Value = 11959.125  # L130 was 12209.5716
$ws.Cells.Item(130, 13).Value = 2020  # M130 was 1648
$ws.Cells.Item(130, 14).Value = -21999.125  # N130 was -22249.5716
$ws.Cells.Item(131, 8).Value = 10870551  # H131 was 10417638
$ws.Cells.Item(131, 9).Value = 125002740  # I131 was 100002270
$ws.Cells.Item(131, 10).Value = 818.7143  # J131 was 820.3953
$ws.Cells.Item(131, 11).Value = 375008220  # K131 was 300006810
$ws.Cells.Item(131, 12).Value = 2456.1429  # L131 was 2461.1859
$ws.Cells.Item(131, 13).Value = -375003180  # M131 was -300001770
$ws.Cells.Item(131, 14).Value = -12536.1429  # N131 was -12541.1859
$ws.Cells.Item(137, 8).Value = 3595  # H137 was 4368.5713
$ws.Cells.Item(137, 9).Value = 2718.4285  # I137 was 5000
$ws.Cells.Item(137, 10).Value = 4152.8184  # J137 was 4196.364
$ws.Cells.Item(137, 11).Value = 8155.2855  # K137 was 15000
$ws.Cells.Item(137, 12).Value = 12458.4552  # L137 was 12589.092
$ws.Cells.Item(137, 13).Value = -3055.2855  # M137 was -9900
$ws.Cells.Item(137, 14).Value = -22658.4552  # N137 was -22789.092
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 6162.8  # H70 was 6305.3687
$ws.Cells.Item(70, 9).Value = 5824.069  # I70 was 5999.6294
$ws.Cells.Item(70, 11).Value = 5824.069  # K70 was 5999.6294
$ws.Cells.Item(70, 13).Value = -5554.069  # M70 was -5729.6294
$ws.Cells.Item(73, 8).Value = 6162.8  # H73 was 6305.3687
$ws.Cells.Item(73, 9).Value = 5824.069  # I73 was 5999.6294
$ws.Cells.Item(73, 11).Value = 5824.069  # K73 was 5999.6294
$ws.Cells.Item(73, 13).Value = -4888.069  # M73 was -5063.6294
$ws.Cells.Item(136, 8).Value = 14149.458  # H136 was 17081.5
$ws.Cells.Item(136, 10).Value = 14149.458  # J136 was 17081.5
$ws.Cells.Item(136, 12).Value = 42448.374  # L136 was 51244.5
$ws.Cells.Item(136, 14).Value = -47548.374  # N136 was -56344.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 2038.4615  # H16 was 2308.9092
$ws.Cells.Item(16, 9).Value = 1966.6666  # I16 was 2371.1428
$ws.Cells.Item(16, 11).Value = 1966.6666  # K16 was 2371.1428
$ws.Cells.Item(16, 13).Value = -1796.6666  # M16 was -2201.1428
$ws.Cells.Item(122, 8).Value = 4054.7742  # H122 was 3624.6758
$ws.Cells.Item(122, 9).Value = 3251.7778  # I122 was 2969.1875
$ws.Cells.Item(122, 10).Value = 9475  # J122 was 7819.8
$ws.Cells.Item(122, 11).Value = 9755.3334  # K122 was 8907.5625
$ws.Cells.Item(122, 12).Value = 28425  # L122 was 23459.4
$ws.Cells.Item(122, 13).Value = -7305.3334  # M122 was -6457.5625
$ws.Cells.Item(122, 14).Value = -33325  # N122 was -28359.4
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 8134546  # H132 was 8776686
$ws.Cells.Item(132, 9).Value = 7009.125  # I132 was 8449.691999999999
$ws.Cells.Item(132, 11).Value = 21027.375  # K132 was 25349.076
$ws.Cells.Item(132, 13).Value = -18497.375  # M132 was -22819.076
